# Update Mental Health Ontology mapping to LSRs.xlsx
# Replace all occurrences of "measurement datum at followup" with
# "measurement datum at post-intervention" in column H ("Class label").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "measurement datum at followup"
$newText = "measurement datum at post-intervention"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H
    if ($cell.Text -eq $oldText) {
        $cell.Value = $newText
    }
}
